# "bug with empty notes is resolved"
#
# Footnote 22 is an empty placeholder note - its body carries no real
# content, just the closing Tibetan punctuation mark - so remove it
# entirely (both its in-text reference mark and its note body).
# Footnote 21's note text had picked up a stray trailing "a" typo;
# strip it so the note reads correctly.

$d = $word.ActiveDocument

# Fix the typo'd trailing "a" in footnote 21's text.
$fn21 = $d.Footnotes.Item(1)
$text21 = $fn21.Range.Text
if ($text21.Length -gt 0 -and $text21.Substring($text21.Length - 1, 1) -eq "a") {
    $fn21.Range.Text = $text21.Substring(0, $text21.Length - 1)
}

# Drop the empty footnote 22 - removes the reference mark in the body
# text and its (empty) note in one call.
$fn22 = $d.Footnotes.Item(2)
$fn22.Delete()
